$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H59").Value = 14500
$ws.Range("J59").Value = 18666.666
$ws.Range("L59").Value = 55999.99800000001
$ws.Range("N59").Value = -57113.99800000001

$ws.Range("H64").Value = 14511.333
$ws.Range("I64").Value = 23040.4
$ws.Range("J64").Value = 3850
$ws.Range("K64").Value = 23040.4
$ws.Range("L64").Value = 3850
$ws.Range("M64").Value = -22792.4
$ws.Range("N64").Value = -4346

$ws.Range("H67").Value = 14511.333
$ws.Range("I67").Value = 23040.4
$ws.Range("J67").Value = 3850
$ws.Range("K67").Value = 23040.4
$ws.Range("L67").Value = 3850
$ws.Range("M67").Value = -22182.4
$ws.Range("N67").Value = -5566

$ws.Range("H74").Value = 6226.636
$ws.Range("I74").Value = 8198.6
$ws.Range("J74").Value = 4583.3335
$ws.Range("K74").Value = 8198.6
$ws.Range("L74").Value = 4583.3335
$ws.Range("M74").Value = -7262.6
$ws.Range("N74").Value = -6455.3335

$ws.Range("H76").Value = 45836170
$ws.Range("I76").Value = 50002730
$ws.Range("K76").Value = 50002730
$ws.Range("M76").Value = -50002415

$ws.Range("H77").Value = 6226.636
$ws.Range("I77").Value = 8198.6
$ws.Range("J77").Value = 4583.3335
$ws.Range("K77").Value = 40993
$ws.Range("L77").Value = 22916.6675
$ws.Range("M77").Value = -36313
$ws.Range("N77").Value = -32276.6675

$ws.Range("H79").Value = 45836170
$ws.Range("I79").Value = 50002730
$ws.Range("K79").Value = 50002730
$ws.Range("M79").Value = -50001638

$ws.Range("H92").Value = 3606.5
$ws.Range("I92").Value = 4154.8
$ws.Range("K92").Value = 4154.8
$ws.Range("M92").Value = -2906.8


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 5250
$ws.Range("J50").Value = 5250
$ws.Range("L50").Value = 5250
$ws.Range("N50").Value = -6678

$ws.Range("H61").Value = 1716
$ws.Range("I61").Value = 1180.5
$ws.Range("K61").Value = 1180.5
$ws.Range("M61").Value = -968.5

$ws.Range("H122").Value = 1624.8
$ws.Range("I122").Value = 1624.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4874.4
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws.Range("H136").Value = 1716
$ws.Range("I136").Value = 1180.5
$ws.Range("K136").Value = 3541.5
$ws.Range("M136").Value = -991.5


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 18239.666

$ws.Range("H85").Value = 18239.666


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 7819.7
$ws.Range("I62").Value = 8516.25
$ws.Range("J62").Value = 6774.875
$ws.Range("K62").Value = 8516.25
$ws.Range("L62").Value = 6774.875
$ws.Range("M62").Value = -7892.25
$ws.Range("N62").Value = -8022.875

$ws.Range("H65").Value = 7819.7
$ws.Range("I65").Value = 8516.25
$ws.Range("J65").Value = 6774.875
$ws.Range("K65").Value = 42581.25
$ws.Range("L65").Value = 33874.375
$ws.Range("M65").Value = -39461.25
$ws.Range("N65").Value = -40114.375

$ws.Range("H122").Value = 1112706.9
$ws.Range("I122").Value = 1251582.8
$ws.Range("J122").Value = 1700
$ws.Range("K122").Value = 3754748.4
$ws.Range("L122").Value = 5100
$ws.Range("M122").Value = -3752298.4
$ws.Range("N122").Value = -10000

$ws.Range("H132").Value = 2411.5
$ws.Range("I132").Value = 1561.3334
$ws.Range("J132").Value = 4962
$ws.Range("K132").Value = 4684.0002
$ws.Range("L132").Value = 14886
$ws.Range("M132").Value = -2154.0002
$ws.Range("N132").Value = -19946

$ws.Range("H140").Value = 51760
$ws.Range("J140").Value = 51760
$ws.Range("L140").Value = 51760
$ws.Range("N140").Value = -62120


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").ClearContents()
$ws.Range("M36").ClearContents()
$ws.Range("N36").ClearContents()

$ws.Range("H42").Value = 1842.25
$ws.Range("J42").Value = 2289.6667
$ws.Range("L42").Value = 6869.000100000001
$ws.Range("N42").Value = -7937.000100000001

$ws.Range("H97").Value = 1566.6666
$ws.Range("J97").Value = 1566.6666
$ws.Range("L97").Value = 4699.9998
$ws.Range("N97").Value = -5691.9998


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4336.364
$ws.Range("I70").Value = 4100
$ws.Range("J70").Value = 4966.6665
$ws.Range("K70").Value = 4100
$ws.Range("L70").Value = 4966.6665
$ws.Range("M70").Value = -3830
$ws.Range("N70").Value = -5506.6665

$ws.Range("H73").Value = 4336.364
$ws.Range("I73").Value = 4100
$ws.Range("J73").Value = 4966.6665
$ws.Range("K73").Value = 4100
$ws.Range("L73").Value = 4966.6665
$ws.Range("M73").Value = -3164
$ws.Range("N73").Value = -6838.6665

$ws.Range("H122").Value = 2168.7
$ws.Range("I122").Value = 2187.4443
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 6562.3329
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -4112.3329
$ws.Range("N122").Value = -10900


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2285.5715
$ws.Range("I16").Value = 2285.5715
$ws.Range("K16").Value = 2285.5715
$ws.Range("M16").Value = -2115.5715

$ws.Range("H31").Value = 638.3333
$ws.Range("I31").Value = 707.5
$ws.Range("J31").Value = 500
$ws.Range("K31").Value = 707.5
$ws.Range("L31").Value = 500
$ws.Range("M31").Value = -459.5
$ws.Range("N31").Value = -996

$ws.Range("H99").Value = 16500
$ws.Range("J99").Value = 16500
$ws.Range("L99").Value = 16500
$ws.Range("N99").Value = -22490

$ws.Range("H136").Value = 3276.5454
$ws.Range("I136").Value = 1457.6666
$ws.Range("J136").Value = 5459.2
$ws.Range("K136").Value = 4372.9998
$ws.Range("L136").Value = 16377.6
$ws.Range("M136").Value = -1822.9998
$ws.Range("N136").Value = -21477.6

$ws.Range("H139").Value = 35350
$ws.Range("J139").Value = 35350
$ws.Range("L139").Value = 35350
$ws.Range("N139").Value = -45630


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 673.25
$ws.Range("I126").Value = 602
$ws.Range("J126").Value = 697
$ws.Range("K126").Value = 1806
$ws.Range("L126").Value = 2091
$ws.Range("M126").Value = 664
$ws.Range("N126").Value = -7031

$ws.Range("H136").Value = 2384.3877
$ws.Range("I136").Value = 2872.742
$ws.Range("J136").Value = 1543.3334
$ws.Range("K136").Value = 8618.226000000001
$ws.Range("L136").Value = 4630.0002
$ws.Range("M136").Value = -6068.226000000001
$ws.Range("N136").Value = -9730.0002

